$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers: BTec_Logo-Orange inline pictures rename image2.jpg -> image1.jpg ---
for ($h = 1; $h -le 2; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
            $shp = $hdr.Range.InlineShapes.Item($i)
            $shp.Name = "image1.jpg"
        }
    }
}

# --- Footers: PearsonLogo.png inline pictures rename image1.png -> image2.png ---
# Direct property assignment on footer-derived InlineShape objects does not
# stick in this host, so route the write through the Selection object.
for ($f = 1; $f -le 2; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
            $shp = $ftr.Range.InlineShapes.Item($i)
            [void]$shp.Select()
            $word.Selection.InlineShapes.Item(1).Name = "image2.png"
        }
    }
}
